$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D, rows 2-7: previously shared formulas referencing row+6/D13,
#     now recorded as fixed observed values (no formula) -------------------
$ws.Range("D2").Value = 45
$ws.Range("D3").Value = 44
$ws.Range("D4").Value = 43
$ws.Range("D5").Value = 43
$ws.Range("D6").Value = 42
$ws.Range("D7").Value = 48

# --- Column D, rows 8-13: now "= (row+6) + 2" -------------------------------
$ws.Range("D8:D11").Formula = "=D14+2"
$ws.Range("D12").Formula = "=D18+2"
$ws.Range("D13").Formula = "=D19+2"

# --- Column D, rows 14-19: now "= (row+6) + 8" ------------------------------
$ws.Range("D14:D18").Formula = "=D20+8"
$ws.Range("D19").Formula = "=D25+8"

# --- Column D, rows 20-25: new raw/physical sampling data (plain values) ---
$ws.Range("D20").Value = 35
$ws.Range("D21").Value = 34
$ws.Range("D22").Value = 33
$ws.Range("D23").Value = 33
$ws.Range("D24").Value = 32
$ws.Range("D25").Value = 38

# --- Selection moves to L19 -------------------------------------------------
$ws.Range("L19").Select()

$wb.Save()
